# The document ends with an (otherwise empty) paragraph that only holds
# the Word-managed "_GoBack" bookmark (left over from the last edit
# position). The commit adds a single run of text ("m") inside that
# paragraph, placed before the bookmark.
$d = $word.ActiveDocument

# Locate the trailing "_GoBack" bookmark. It is addressable by name even
# though Word hides it from the default Bookmarks enumeration/Count, so
# this is more robust than hard-coding a character offset.
$bm = $d.Bookmarks("_GoBack")

# A zero-length range right at the bookmark's start is exactly the
# insertion point used when Word types new text just before it.
$insertionPoint = $d.Range($bm.Start, $bm.Start)
$insertionPoint.InsertBefore("m")

Write-Output "Inserted 'm' before _GoBack bookmark."
